# Updated cryptos list with latest prices / volume data, and swap ARBITRUM / RenderToken row order
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    'D2' = '26.521.44'
    'E2' = '  -2.56%  '
    'D3' = '1.812.32'
    'E3' = '  -2.23%  '
    'D4' = '1.007'
    'E4' = '  +0.64%  '
    'E5' = '  +0.60%  '
    'D6' = '308.32'
    'E6' = '  -1.74%  '
    'D7' = '0.4562'
    'E7' = '  -1.91%  '
    'D8' = '0.3666'
    'E8' = '  -1.23%  '
    'D9' = '0.07133'
    'E9' = '  -2.22%  '
    'D10' = '0.8796'
    'E10' = '  -1.23%  '
    'E11' = '  -1.36%  '
    'D12' = '19.38'
    'E12' = '  -3.57%  '
    'D13' = '1.802.67'
    'E13' = '  -0.44%  '
    'D14' = '5.290'
    'E14' = '  -2.09%  '
    'D15' = '6.374'
    'E15' = '  -2.23%  '
    'D16' = '86.63'
    'E16' = '  -5.08%  '
    'D18' = '0.000008588'
    'E18' = '  -3.64%  '
    'D20' = '26.589.61'
    'E20' = '  -2.41%  '
    'D21' = '14.24'
    'E21' = '  -3.16%  '
    'D22' = '5.012'
    'E22' = '  -1.51%  '
    'E23' = '  -0.53%  '
    'D24' = '1.990'
    'E24' = '  +1.57%  '
    'D25' = '151.55'
    'E25' = '  +0.07%  '
    'D26' = '17.96'
    'E26' = '  -2.35%  '
    'D27' = '2.063'
    'E27' = '  +1.11%  '
    'D28' = '112.82'
    'E28' = '  -2.62%  '
    'E29' = '  -3.94%  '
    'D30' = '0.08688'
    'E30' = '  -1.67%  '
    'D31' = '3.061'
    'E31' = '  -2.50%  '
    'D32' = '4.522'
    'E32' = '  +0.10%  '
    'D33' = '0.7323'
    'E33' = '  -5.18%  '
    'B34' = 'ARBITRUM'
    'C34' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D34' = '1.119'
    'E34' = '  -4.10%  '
    'B35' = 'RenderToken'
    'C35' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D35' = '2.673'
    'E35' = '  -1.45%  '
    'E36' = '  +1.60%  '
    'E37' = '  -2.23%  '
    'D38' = '0.01956'
    'E38' = '  +0.69%  '
    'D39' = '0.05111'
    'E39' = '  -2.07%  '
    'D40' = '2.893'
    'E40' = '  -1.83%  '
    'D41' = '6.982'
    'E41' = '  -1.02%  '
    'E42' = '  -2.43%  '
    'D43' = '0.1561'
    'E43' = '  -3.92%  '
    'E44' = '  -3.66%  '
    'D45' = '1.008'
    'E45' = '  +0.70%  '
    'D46' = '0.4607'
    'E46' = '  -3.81%  '
    'E47' = '  -3.31%  '
    'D48' = '100.97'
    'E48' = '  -1.65%  '
    'E49' = '  -3.19%  '
    'D50' = '0.05998'
    'E50' = '  -3.23%  '
    'D51' = '64.41'
    'E51' = '  -1.42%  '
}

foreach ($cellRef in $updates.Keys) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $updates[$cellRef]
    $c.Style = "Normal"
}
